$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 2.63
$ws.Range("K2").Value = 2.38
$ws.Range("N2").Value = 15
$ws.Range("Y2").Value = 1.53
$ws.Range("Z2").Value = 2.38
$ws.Range("AH3").Value = 7
$ws.Range("AI3").Value = 15
$ws.Range("AK3").Value = 251
$ws.Range("AO3").Value = 34
$ws.Range("K3").Value = 2.2
$ws.Range("Q3").Value = 1.91
$ws.Range("R3").Value = 1.99
$ws.Range("U3").Value = 3.2
$ws.Range("V3").Value = 1.36
$ws.Range("Y3").Value = 1.7
$ws.Range("Z3").Value = 2.05
$ws.Range("AM4").Value = 15
$ws.Range("G4").Value = 2.7
$ws.Range("I4").Value = 2.55
$ws.Range("AG5").Value = 29
$ws.Range("AH5").Value = 17
$ws.Range("AP5").Value = 81
$ws.Range("G5").Value = 1.2
$ws.Range("H5").Value = 7.5
$ws.Range("I5").Value = 13
$ws.Range("J5").Value = 1.53
$ws.Range("K5").Value = 3.2
$ws.Range("L5").Value = 9
$ws.Range("N5").Value = 26
$ws.Range("U5").Value = 1.73
$ws.Range("V5").Value = 2.1
$ws.Range("W5").Value = 1.17
$ws.Range("X5").Value = 5
$ws.Range("Y5").Value = 1.75
$ws.Range("Z5").Value = 2
$ws.Range("AA7").Value = 5
$ws.Range("AB7").Value = 7
$ws.Range("AD7").Value = 13
$ws.Range("AL7").Value = 10
$ws.Range("AM7").Value = 23
$ws.Range("AN7").Value = 19
$ws.Range("AP7").Value = 51
$ws.Range("AR7").Value = 1.95
$ws.Range("AS7").Value = 1.9
$ws.Range("G7").Value = 1.8
$ws.Range("I7").Value = 5.25
$ws.Range("J7").Value = 2.5
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 7
$ws.Range("AA8").Value = 8
$ws.Range("AD8").Value = 26
$ws.Range("AF8").Value = 34
$ws.Range("AG8").Value = 8.5
$ws.Range("AM8").Value = 13
$ws.Range("K8").Value = 2.05
$ws.Range("N8").Value = 8.5
$ws.Range("Y8").Value = 1.8
$ws.Range("Z8").Value = 1.91
$ws.Range("AG11").Value = 11
$ws.Range("AH11").Value = 7
$ws.Range("AM11").Value = 10
$ws.Range("AQ11").Value = 26
$ws.Range("I11").Value = 2
$ws.Range("K11").Value = 2.2
$ws.Range("O11").Value = 1.29
$ws.Range("P11").Value = 3.5
$ws.Range("Q11").Value = 1.98
$ws.Range("R11").Value = 1.88
$ws.Range("AN12").Value = 11
$ws.Range("AR12").Value = 1.88
$ws.Range("AS12").Value = 1.98
$ws.Range("H12").Value = 3.1
$ws.Range("K12").Value = 1.91
$ws.Range("W12").Value = 1.57
$ws.Range("X12").Value = 2.25
$ws.Range("AA13").Value = 9.5
$ws.Range("AB13").Value = 11
$ws.Range("AE13").Value = 13
$ws.Range("AF13").Value = 21
$ws.Range("AG13").Value = 15
$ws.Range("AH13").Value = 7.5
$ws.Range("AK13").Value = 126
$ws.Range("AL13").Value = 15
$ws.Range("AP13").Value = 26
$ws.Range("AQ13").Value = 29
$ws.Range("G13").Value = 1.83
$ws.Range("H13").Value = 3.7
$ws.Range("J13").Value = 2.4
$ws.Range("K13").Value = 2.38
$ws.Range("M13").Value = 1.03
$ws.Range("N13").Value = 15
$ws.Range("O13").Value = 1.18
$ws.Range("P13").Value = 4.5
$ws.Range("Q13").Value = 1.62
$ws.Range("R13").Value = 2.25
$ws.Range("U13").Value = 2.5
$ws.Range("V13").Value = 1.5
$ws.Range("W13").Value = 1.3
$ws.Range("X13").Value = 3.4
$ws.Range("Y13").Value = 1.57
$ws.Range("Z13").Value = 2.25
$ws.Range("AB14").Value = 7
$ws.Range("AD14").Value = 9.5
$ws.Range("AL14").Value = 21
$ws.Range("AQ14").Value = 51
$ws.Range("G14").Value = 1.4
$ws.Range("I14").Value = 7.5
$ws.Range("Y14").Value = 1.91
$ws.Range("Z14").Value = 1.91
$ws.Range("AB15").Value = 21
$ws.Range("AC15").Value = 13
$ws.Range("AE15").Value = 29
$ws.Range("AK15").Value = 151
$ws.Range("AL15").Value = 8.5
$ws.Range("G15").Value = 3.9
$ws.Range("I15").Value = 1.75
$ws.Range("G16").Value = 1.13
$ws.Range("I17").Value = 1.42
$ws.Range("O17").Value = 1.14
$ws.Range("P17").Value = 5
$ws.Range("S17").Value = 1.78
$ws.Range("T17").Value = 2.03
$ws.Range("U17").Value = 2.2
$ws.Range("V17").Value = 1.62
$ws.Range("AA18").Value = 9
$ws.Range("AF18").Value = 26
$ws.Range("AG18").Value = 12
$ws.Range("AH18").Value = 6.5
$ws.Range("AL18").Value = 11
$ws.Range("AN18").Value = 12
$ws.Range("AP18").Value = 23
$ws.Range("G18").Value = 2.25
$ws.Range("H18").Value = 3.4
$ws.Range("I18").Value = 2.88
$ws.Range("J18").Value = 2.88
$ws.Range("K18").Value = 2.2
$ws.Range("M18").Value = 1.04
$ws.Range("N18").Value = 9
$ws.Range("O18").Value = 1.25
$ws.Range("P18").Value = 3.75
$ws.Range("Q18").Value = 1.8
$ws.Range("R18").Value = 2
$ws.Range("U18").Value = 3
$ws.Range("V18").Value = 1.36
$ws.Range("W18").Value = 1.36
$ws.Range("X18").Value = 3
$ws.Range("Y18").Value = 1.67
$ws.Range("Z18").Value = 2.1
$ws.Range("AC19").Value = 11
$ws.Range("AF19").Value = 26
$ws.Range("AH19").Value = 17
$ws.Range("AJ19").Value = 67
$ws.Range("G19").Value = 1.17
$ws.Range("J19").Value = 1.5
$ws.Range("L19").Value = 10
$ws.Range("Q19").Value = 1.33
$ws.Range("O20").Value = 1.29
$ws.Range("P20").Value = 3.5
$ws.Range("Q20").Value = 1.93
$ws.Range("R20").Value = 1.93
$ws.Range("U20").Value = 3.25
$ws.Range("V20").Value = 1.33
$ws.Range("AD21").Value = 29
$ws.Range("AE21").Value = 23
$ws.Range("AG21").Value = 8.5
$ws.Range("AI21").Value = 15
$ws.Range("AJ21").Value = 51
$ws.Range("AL21").Value = 8
$ws.Range("AM21").Value = 12
$ws.Range("AN21").Value = 10
$ws.Range("AO21").Value = 23
$ws.Range("AQ21").Value = 34
$ws.Range("G21").Value = 2.88
$ws.Range("I21").Value = 2.5
$ws.Range("J21").Value = 3.4
$ws.Range("M21").Value = 1.05
$ws.Range("N21").Value = 8.5
$ws.Range("O21").Value = 1.33
$ws.Range("P21").Value = 3.25
$ws.Range("Q21").Value = 2.08
$ws.Range("R21").Value = 1.73
$ws.Range("U21").Value = 3.75
$ws.Range("V21").Value = 1.25
$ws.Range("Y21").Value = 1.83
$ws.Range("Z21").Value = 1.83
$ws.Range("AB22").Value = 6
$ws.Range("AG22").Value = 10
$ws.Range("AH22").Value = 9
$ws.Range("H22").Value = 4.5
$ws.Range("I22").Value = 6.25
$ws.Range("N22").Value = 10
$ws.Range("Y22").Value = 2.25
$ws.Range("Z22").Value = 1.57
